$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 344.16666
$ws.Range("I53").Value = 312.25
$ws.Range("J53").Value = 360.125
$ws.Range("K53").Value = 312.25
$ws.Range("L53").Value = 360.125
$ws.Range("M53").Value = 324.75
$ws.Range("N53").Value = -1634.125
$ws.Range("H62").Value = 3502.25
$ws.Range("I62").Value = 2004.5
$ws.Range("K62").Value = 2004.5
$ws.Range("M62").Value = -1380.5
$ws.Range("H65").Value = 3502.25
$ws.Range("I65").Value = 2004.5
$ws.Range("K65").Value = 10022.5
$ws.Range("M65").Value = -6902.5
$ws.Range("H76").Value = 4634
$ws.Range("I76").Value = 4751.7144
$ws.Range("K76").Value = 4751.7144
$ws.Range("M76").Value = -4436.7144
$ws.Range("H79").Value = 4634
$ws.Range("I79").Value = 4751.7144
$ws.Range("K79").Value = 4751.7144
$ws.Range("M79").Value = -3659.7144
$ws.Range("H107").Value = 380.5
$ws.Range("I107").Value = 71
$ws.Range("K107").Value = 71
$ws.Range("M107").Value = 1849
$ws.Range("H132").Value = 1652.8518
$ws.Range("I132").Value = 1331.8462
$ws.Range("K132").Value = 3995.5386
$ws.Range("M132").Value = -1465.5386
$ws.Range("H137").Value = 2074.875
$ws.Range("I137").Value = 2014.1428
$ws.Range("K137").Value = 6042.428400000001
$ws.Range("M137").Value = -3492.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 667.36365
$ws.Range("I2").Value = 346.85715
$ws.Range("K2").Value = 346.85715
$ws.Range("M2").Value = -233.85715
$ws.Range("H32").Value = 2757.516
$ws.Range("I32").Value = 2893.5173
$ws.Range("K32").Value = 2893.5173
$ws.Range("M32").Value = -2606.5173
$ws.Range("H39").Value = 1449.5
$ws.Range("I39").Value = 1449.5
$ws.Range("K39").Value = 1449.5
$ws.Range("M39").Value = -929.5
$ws.Range("H61").Value = 747.5
$ws.Range("I61").Value = 747.5
$ws.Range("K61").Value = 747.5
$ws.Range("M61").Value = -535.5
$ws.Range("H74").Value = 949.1111
$ws.Range("I74").Value = 949.1111
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 949.1111
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -75.11109999999996
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 949.1111
$ws.Range("I77").Value = 949.1111
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 4745.555499999999
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -377.5554999999995
$ws.Range("N77").ClearContents()
$ws.Range("H116").Value = 667.36365
$ws.Range("I116").Value = 346.85715
$ws.Range("K116").Value = 346.85715
$ws.Range("M116").Value = 1947.14285
$ws.Range("H122").Value = 1914.7142
$ws.Range("I122").Value = 1000.2857
$ws.Range("K122").Value = 3000.8571
$ws.Range("M122").Value = -550.8571000000002
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 747.5
$ws.Range("I136").Value = 747.5
$ws.Range("K136").Value = 2242.5
$ws.Range("M136").Value = 307.5
$ws.Range("H139").Value = 80650
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 667.36365
$ws.Range("I3").Value = 346.85715
$ws.Range("K3").Value = 346.85715
$ws.Range("M3").Value = -232.85715
$ws.Range("H33").Value = 12873.667
$ws.Range("I33").Value = 5560.5
$ws.Range("K33").Value = 5560.5
$ws.Range("M33").Value = -5224.5
$ws.Range("H134").Value = 2657.261
$ws.Range("I134").Value = 2709.9546
$ws.Range("K134").Value = 8129.8638
$ws.Range("M134").Value = -5594.8638
$ws.Range("H135").Value = 45000
$ws.Range("J135").Value = 45000
$ws.Range("L135").Value = 45000
$ws.Range("N135").Value = -55140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3152.3333
$ws.Range("I16").Value = 2874.8
$ws.Range("K16").Value = 2874.8
$ws.Range("M16").Value = -2587.8
$ws.Range("H31").Value = 2460.3333
$ws.Range("I31").Value = 2817.8572
$ws.Range("K31").Value = 2817.8572
$ws.Range("M31").Value = -2522.8572
$ws.Range("H34").Value = 2460.3333
$ws.Range("I34").Value = 2817.8572
$ws.Range("K34").Value = 2817.8572
$ws.Range("M34").Value = -2615.8572
$ws.Range("H35").Value = 182.25
$ws.Range("I35").Value = 182.25
$ws.Range("K35").Value = 182.25
$ws.Range("M35").Value = 111.75
$ws.Range("H105").Value = 2426.5
$ws.Range("I105").Value = 1800.4
$ws.Range("K105").Value = 1800.4
$ws.Range("M105").Value = -53.40000000000009
$ws.Range("H107").Value = 1372.8572
$ws.Range("I107").Value = 1199.625
$ws.Range("K107").Value = 1199.625
$ws.Range("M107").Value = 720.375
$ws.Range("H113").Value = 3152.3333
$ws.Range("I113").Value = 2874.8
$ws.Range("K113").Value = 2874.8
$ws.Range("M113").Value = -704.8000000000002
$ws.Range("H122").Value = 1779.8182
$ws.Range("I122").Value = 842.1111
$ws.Range("K122").Value = 2526.3333
$ws.Range("M122").Value = -76.33329999999978
$ws.Range("H134").Value = 906.5294
$ws.Range("I134").Value = 906.5294
$ws.Range("K134").Value = 2719.5882
$ws.Range("M134").Value = -184.5882000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 259.75
$ws.Range("I10").Value = 329.66666
$ws.Range("K10").Value = 988.9999799999999
$ws.Range("M10").Value = -849.9999799999999
$ws.Range("H46").Value = 1766.6666
$ws.Range("I46").Value = 1766.6666
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 5299.9998
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -5208.9998
$ws.Range("N46").ClearContents()
$ws.Range("H128").Value = 636159.2
$ws.Range("I128").Value = 636159.2
$ws.Range("K128").Value = 1908477.6
$ws.Range("M128").Value = -1903497.6
$ws.Range("H140").Value = 8206.666999999999
$ws.Range("I140").Value = 1012.5
$ws.Range("J140").Value = 17199.375
$ws.Range("K140").Value = 3037.5
$ws.Range("L140").Value = 51598.125
$ws.Range("M140").Value = 2142.5
$ws.Range("N140").Value = -61958.125
$ws.Range("H141").Value = 10645.75
$ws.Range("I141").Value = 10645.75
$ws.Range("K141").Value = 31937.25
$ws.Range("M141").Value = -26757.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2122
$ws.Range("I80").Value = 2088.5
$ws.Range("J80").Value = 2166.6667
$ws.Range("K80").Value = 2088.5
$ws.Range("L80").Value = 2166.6667
$ws.Range("M80").Value = -1090.5
$ws.Range("N80").Value = -4162.6667
$ws.Range("H83").Value = 2122
$ws.Range("I83").Value = 2088.5
$ws.Range("J83").Value = 2166.6667
$ws.Range("K83").Value = 10442.5
$ws.Range("L83").Value = 10833.3335
$ws.Range("M83").Value = -5450.5
$ws.Range("N83").Value = -20817.3335
$ws.Range("H102").Value = 5120.25
$ws.Range("I102").Value = 4993.6665
$ws.Range("K102").Value = 4993.6665
$ws.Range("M102").Value = -3371.6665
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 848.1429000000001
$ws.Range("I16").Value = 848.1429000000001
$ws.Range("K16").Value = 848.1429000000001
$ws.Range("M16").Value = -678.1429000000001
$ws.Range("H132").Value = 4191.7827
$ws.Range("I132").Value = 4433.8887
$ws.Range("K132").Value = 13301.6661
$ws.Range("M132").Value = -10771.6661
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 42499.5
$ws.Range("J54").Value = 44999
$ws.Range("L54").Value = 44999
$ws.Range("N54").Value = -46039
$ws.Range("H55").Value = 20525
$ws.Range("I55").Value = 11050
$ws.Range("K55").Value = 11050
$ws.Range("M55").Value = -10773
$ws.Range("H58").Value = 5500
$ws.Range("I58").Value = 4000
$ws.Range("K58").Value = 4000
$ws.Range("M58").Value = -3692
$ws.Range("H62").Value = 12533.556
$ws.Range("I62").Value = 22901
$ws.Range("J62").Value = 9571.429
$ws.Range("K62").Value = 22901
$ws.Range("L62").Value = 9571.429
$ws.Range("M62").Value = -22277
$ws.Range("N62").Value = -10819.429
$ws.Range("H65").Value = 12533.556
$ws.Range("I65").Value = 22901
$ws.Range("J65").Value = 9571.429
$ws.Range("K65").Value = 114505
$ws.Range("L65").Value = 47857.145
$ws.Range("M65").Value = -111385
$ws.Range("N65").Value = -54097.145
$ws.Range("H81").Value = 912598.5600000001
$ws.Range("I81").Value = 1010.375
$ws.Range("K81").Value = 2020.75
$ws.Range("M81").Value = -959.75
$ws.Range("H84").Value = 912598.5600000001
$ws.Range("I84").Value = 1010.375
$ws.Range("K84").Value = 10103.75
$ws.Range("M84").Value = -4799.75
$ws.Range("H96").Value = 969.75
$ws.Range("I96").Value = 947.5
$ws.Range("J96").Value = 992
$ws.Range("K96").Value = 947.5
$ws.Range("L96").Value = 992
$ws.Range("M96").Value = 425.5
$ws.Range("N96").Value = -3738
